# Add team record (Wins/Losses/Ties) columns to the MIL_1998 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold, centered/top-aligned, thin border)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill in the team record for every player row (2-45): 74 wins, 88 losses, 0 ties
$ws.Range("AD2:AD45").Value = 74
$ws.Range("AE2:AE45").Value = 88
$ws.Range("AF2:AF45").Value = 0
